# Treasure_coord_gold.xlsx edit: updated pickup icons and fuel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the D1 cell comment text (pickup-type legend) ---
# Was:
#   Author:
#   1. coin
#   2. small chest
#   3. large chest
# Now:
#   Author:
#   1. coin
#   2. chest
#   3. fuel
$ws.Range("D1").Comment.Text("Author:`r`n1. coin`r`n2. chest`r`n3. fuel")

# --- Fill in the "type" (D) column for rows that were missing it, and
#     bump a handful of rows to the new "fuel" type (3), maxing out
#     their amount columns at 500 to match the new fuel pickup values ---
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 500
$ws.Range("F4").Value = 500
$ws.Range("G4").Value = 500

$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1

$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 500
$ws.Range("F8").Value = 500
$ws.Range("G8").Value = 500

$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("D11").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1

$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 500
$ws.Range("F15").Value = 500
$ws.Range("G15").Value = 500

$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("D23").Value = 2

$ws.Range("D24").Value = 3
$ws.Range("E24").Value = 500
$ws.Range("F24").Value = 500
$ws.Range("G24").Value = 500

$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 430
$ws.Range("F25").Value = 100
$ws.Range("G25").Value = 30

$ws.Range("D26").Value = 2

$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 500
$ws.Range("F27").Value = 500
$ws.Range("G27").Value = 500

$ws.Range("D28").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("D30").Value = 1
$ws.Range("D31").Value = 2
$ws.Range("D32").Value = 1

$ws.Range("D33").Value = 3
$ws.Range("F33").Value = 500
$ws.Range("G33").Value = 500

$ws.Range("D34").Value = 2
$ws.Range("D35").Value = 2
$ws.Range("D36").Value = 2
$ws.Range("D37").Value = 2
$ws.Range("D38").Value = 2
$ws.Range("D39").Value = 3
$ws.Range("D40").Value = 2
$ws.Range("D41").Value = 2
$ws.Range("D42").Value = 2

$ws.Range("D43").Value = 3
$ws.Range("E43").Value = 500
$ws.Range("F43").Value = 500
$ws.Range("G43").Value = 500

$ws.Range("D44").Value = 2
$ws.Range("D45").Value = 2
$ws.Range("D46").Value = 2
$ws.Range("D47").Value = 2
$ws.Range("D48").Value = 2

$ws.Range("D49").Value = 2
$ws.Range("E49").Value = 1000
$ws.Range("F49").Value = 500
$ws.Range("G49").Value = 300

# --- Move the active selection from G15 to D1 ---
$ws.Range("D1").Select()
